$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header for column B (shared string "value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Apply the same date-cell formatting (style) used by A2 to the new rows A3:A22
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Dates (column A) and values (column B) for rows 2..22
$dates = @(38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022)
$values = @($null,1.111105389870159,-0.4885592833739349,-0.2059563123693375,0.1984123724363851,-0.4432005650260806,1.133846722438525,0.6449669885999487,1.147096153021487,1.287777024550762,2.027763112344405,1.45091979290124,1.309848083191967,0.7492024424422983,1.584092467432896,-7.952290978198695,3.004777693925043,4.519371604409206,-0.8677070965151246,0.3575298869986865,$null)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    if ($null -ne $values[$i]) {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    } else {
        $ws.Cells.Item($row, 2).ClearContents()
    }
}
